# Daily attendance processing - 2025-11-04 22:20:38
# Reorders the "Recorded By" (column G) entries so that "System" is moved
# from the front of the comma-separated list to the back (full reverse of
# the list), leaving cells that do not contain an exact "System" entry
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($value -ne $null -and $value -ne "") {
        $parts = $value -split ", "

        if ($parts.Count -gt 1 -and ($parts -contains "System")) {
            $reversedParts = $parts[($parts.Count - 1)..0]
            $cell.Value2 = [string]::Join(", ", $reversedParts)
        }
    }
}
